# Regenerate orders with updated distance/size codes.
#
# The experiment's distance and size condition codes changed:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
# (S20 and S25 are unchanged.)
#
# These codes appear embedded inside many string cells throughout the
# sheet (Condition, Filename_Left, Filename_Right, Distance, Size
# columns, e.g. "Face12_D51_S30" -> "Face12_D55_S31",
# "Fixation_D64_l.png" -> "Fixation_D69_l.png", "D80" -> "D86",
# "S30" -> "S31"), so we walk every used cell and rewrite its text,
# leaving numeric/boolean cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$startRow = $used.Row
$startCol = $used.Column
$numRows = $used.Rows.Count
$numCols = $used.Columns.Count
$endRow = $startRow + $numRows - 1
$endCol = $startCol + $numCols - 1

for ($r = $startRow; $r -le $endRow; $r++) {
    for ($c = $startCol; $c -le $endCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -is [string]) {
            $newVal = $val.Replace("D51", "D55").Replace("D64", "D69").Replace("D80", "D86").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
